{"js": "// Update each \"two-digit x two-digit = product\" answer cell in the table\n// with a freshly generated problem/answer pair, matching the old value to\n// the new value one-for-one (each old string is unique in the document).\n\nconst replacements = [\n  { old: \"92\u00d722=2024\", new: \"48\u00d792=4416\" },\n  { old: \"44\u00d718=792\", new: \"96\u00d741=3936\" },\n  { old: \"80\u00d773=5840\", new: \"50\u00d780=4000\" },\n  { old: \"42\u00d712=504\", new: \"93\u00d730=2790\" },\n  { old: \"55\u00d775=4125\", new: \"34\u00d714=476\" },\n  { old: \"77\u00d766=5082\", new: \"57\u00d745=2565\" },\n  { old: \"66\u00d722=1452\", new: \"22\u00d772=1584\" },\n  { old: \"83\u00d753=4399\", new: \"19\u00d766=1254\" },\n  { old: \"55\u00d764=3520\", new: \"70\u00d742=2940\" },\n  { old: \"96\u00d721=2016\", new: \"73\u00d726=1898\" },\n  { old: \"43\u00d749=2107\", new: \"92\u00d734=3128\" },\n  { old: \"37\u00d767=2479\", new: \"69\u00d722=1518\" },\n  { old: \"92\u00d711=1012\", new: \"32\u00d764=2048\" },\n  { old: \"46\u00d784=3864\", new: \"89\u00d731=2759\" },\n  { old: \"33\u00d753=1749\", new: \"58\u00d793=5394\" },\n  { old: \"35\u00d791=3185\", new: \"52\u00d792=4784\" },\n  { old: \"30\u00d781=2430\", new: \"50\u00d733=1650\" },\n  { old: \"86\u00d772=6192\", new: \"95\u00d718=1710\" },\n  { old: \"96\u00d758=5568\", new: \"54\u00d763=3402\" },\n  { old: \"78\u00d781=6318\", new: \"18\u00d756=1008\" },\n  { old: \"29\u00d797=2813\", new: \"37\u00d743=1591\" },\n  { old: \"99\u00d793=9207\", new: \"99\u00d724=2376\" },\n  { old: \"91\u00d794=8554\", new: \"28\u00d797=2716\" },\n  { old: \"30\u00d797=2910\", new: \"18\u00d784=1512\" },\n  { old: \"84\u00d719=1596\", new: \"50\u00d783=4150\" },\n];\n\nconst body = context.document.body;\n\nfor (const { old, new: replacement } of replacements) {\n  const results = body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each \"two-digit x two-digit = product\" answer cell in the table\n# with a freshly generated problem/answer pair, matching the old value to\n# the new value one-for-one (each old string is unique in the document).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"92\u00d722=2024\"; New = \"48\u00d792=4416\" },\n    @{ Old = \"44\u00d718=792\";  New = \"96\u00d741=3936\" },\n    @{ Old = \"80\u00d773=5840\"; New = \"50\u00d780=4000\" },\n    @{ Old = \"42\u00d712=504\";  New = \"93\u00d730=2790\" },\n    @{ Old = \"55\u00d775=4125\"; New = \"34\u00d714=476\" },\n    @{ Old = \"77\u00d766=5082\"; New = \"57\u00d745=2565\" },\n    @{ Old = \"66\u00d722=1452\"; New = \"22\u00d772=1584\" },\n    @{ Old = \"83\u00d753=4399\"; New = \"19\u00d766=1254\" },\n    @{ Old = \"55\u00d764=3520\"; New = \"70\u00d742=2940\" },\n    @{ Old = \"96\u00d721=2016\"; New = \"73\u00d726=1898\" },\n    @{ Old = \"43\u00d749=2107\"; New = \"92\u00d734=3128\" },\n    @{ Old = \"37\u00d767=2479\"; New = \"69\u00d722=1518\" },\n    @{ Old = \"92\u00d711=1012\"; New = \"32\u00d764=2048\" },\n    @{ Old = \"46\u00d784=3864\"; New = \"89\u00d731=2759\" },\n    @{ Old = \"33\u00d753=1749\"; New = \"58\u00d793=5394\" },\n    @{ Old = \"35\u00d791=3185\"; New = \"52\u00d792=4784\" },\n    @{ Old = \"30\u00d781=2430\"; New = \"50\u00d733=1650\" },\n    @{ Old = \"86\u00d772=6192\"; New = \"95\u00d718=1710\" },\n    @{ Old = \"96\u00d758=5568\"; New = \"54\u00d763=3402\" },\n    @{ Old = \"78\u00d781=6318\"; New = \"18\u00d756=1008\" },\n    @{ Old = \"29\u00d797=2813\"; New = \"37\u00d743=1591\" },\n    @{ Old = \"99\u00d793=9207\"; New = \"99\u00d724=2376\" },\n    @{ Old = \"91\u00d794=8554\"; New = \"28\u00d797=2716\" },\n    @{ Old = \"30\u00d797=2910\"; New = \"18\u00d784=1512\" },\n    @{ Old = \"84\u00d719=1596\"; New = \"50\u00d783=4150\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($null, $true, $true, $false, $null, $null, $true, $null, $null, $r.New, 2) | Out-Null\n}\n"}
